$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$f2 = @"
<rpc-reply message-id="urn:uuid:3f2511dd-86af-495b-b8de-93d20cc711dd">
  <data/>
</rpc-reply>

"@

$ws.Range("F2").Value = $f2

$g2 = @"
  <edit-config>
    <target>
      <candidate/>
    </target>
    <config>
      <network-instances xmlns="http://openconfig.net/yang/network-instance">
        <network-instance>
          <name>Prueba_LxVPN</name>
          <config>
            <name>Prueba_LxVPN</name>
            <type xmlns:oc-ni-types="http://openconfig.net/yang/network-instance-types">oc-ni-types:L3VRF</type>
          </config>
          <protocols>
            <protocol>
              <identifier xmlns:oc-pol-types="http://openconfig.net/yang/policy-types">oc-pol-types:BGP</identifier>
              <name>default</name>
              <config>
                <identifier xmlns:oc-pol-types="http://openconfig.net/yang/policy-types">oc-pol-types:BGP</identifier>
                <name>default</name>
              </config>
              <bgp>
                <global>
                  <config>
                    <as>65000</as>
                  </config>
                </global>
                <neighbors>
                  <neighbor>
                    <neighbor-address>192.168.1.2</neighbor-address>
                    <config>
                      <neighbor-address>192.168.1.2</neighbor-address>
                      <peer-as>65123</peer-as>
                    </config>
                  </neighbor>
                </neighbors>
              </bgp>
            </protocol>
          </protocols>
        </network-instance>
      </network-instances>
    </config>
  </edit-config>
"@

$ws.Range("G2").Value = $g2
